# Updates cryptos.xlsx price/volume figures and swaps the
# Monero / NEARProtocol rows, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a plain text string even when the
    # value looks numeric (e.g. "0.162", "1.00"), so Excel does not
    # silently coerce it into a floating point number.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Row 36 / 37: NEARProtocol and Monero swap places ---
Set-TextValue $ws.Range("B36") "NEARProtocol"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D36") "4.62"
Set-TextValue $ws.Range("E36") "  -3.40%  "

Set-TextValue $ws.Range("B37") "Monero"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D37") "151.69"
Set-TextValue $ws.Range("E37") "  -0.47%  "

# --- Refreshed Price (D) and Volume(1h) (E) values ---
Set-TextValue $ws.Range("D2") "62.727.70"
Set-TextValue $ws.Range("E2") "  -0.49%  "
Set-TextValue $ws.Range("D3") "2.455.98"
Set-TextValue $ws.Range("E3") "  -0.66%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D5") "571.19"
Set-TextValue $ws.Range("E5") "  -1.13%  "
Set-TextValue $ws.Range("D6") "146.04"
Set-TextValue $ws.Range("E6") "  -0.57%  "
Set-TextValue $ws.Range("E7") "  -0.04%  "
Set-TextValue $ws.Range("E8") "  -1.87%  "
Set-TextValue $ws.Range("E9") "  -1.01%  "
Set-TextValue $ws.Range("D10") "0.162"
Set-TextValue $ws.Range("E10") "  -0.38%  "
Set-TextValue $ws.Range("D11") "5.16"
Set-TextValue $ws.Range("E11") "  -2.29%  "
Set-TextValue $ws.Range("E12") "  -1.45%  "
Set-TextValue $ws.Range("E13") "  -1.35%  "
Set-TextValue $ws.Range("E14") "  -3.06%  "
Set-TextValue $ws.Range("D15") "2.895.58"
Set-TextValue $ws.Range("E15") "  -0.83%  "
Set-TextValue $ws.Range("D16") "62.611.82"
Set-TextValue $ws.Range("E16") "  -0.79%  "
Set-TextValue $ws.Range("D17") "2.456.94"
Set-TextValue $ws.Range("E17") "  -0.34%  "
Set-TextValue $ws.Range("D18") "7.69"
Set-TextValue $ws.Range("E18") "  -6.07%  "
Set-TextValue $ws.Range("D20") "2.23"
Set-TextValue $ws.Range("E20") "  +0.22%  "
Set-TextValue $ws.Range("D21") "321.02"
Set-TextValue $ws.Range("E21") "  -2.65%  "
Set-TextValue $ws.Range("E22") "  -0.17%  "
Set-TextValue $ws.Range("E23") "  +0.05%  "
Set-TextValue $ws.Range("D24") "9.87"
Set-TextValue $ws.Range("E24") "  +3.07%  "
Set-TextValue $ws.Range("D25") "64.75"
Set-TextValue $ws.Range("D26") "649.57"
Set-TextValue $ws.Range("E26") "  -2.41%  "
Set-TextValue $ws.Range("D27") "2.579.59"
Set-TextValue $ws.Range("E27") "  -0.47%  "
Set-TextValue $ws.Range("D28") "0.0₃0949"
Set-TextValue $ws.Range("E28") "  -3.85%  "
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  -0.11%  "
Set-TextValue $ws.Range("D30") "1.41"
Set-TextValue $ws.Range("E30") "  -2.63%  "
Set-TextValue $ws.Range("D31") "7.81"
Set-TextValue $ws.Range("E31") "  -3.37%  "
Set-TextValue $ws.Range("E32") "  -3.21%  "
Set-TextValue $ws.Range("E33") "  -0.15%  "
Set-TextValue $ws.Range("E34") "  -0.05%  "
Set-TextValue $ws.Range("E35") "  -3.57%  "
Set-TextValue $ws.Range("E38") "  -1.40%  "
Set-TextValue $ws.Range("E39") "  -2.34%  "
Set-TextValue $ws.Range("D40") "5.29"
Set-TextValue $ws.Range("E40") "  -2.85%  "
Set-TextValue $ws.Range("D41") "2.62"
Set-TextValue $ws.Range("E41") "  -4.13%  "
Set-TextValue $ws.Range("E42") "  -3.97%  "
Set-TextValue $ws.Range("E43") "  +0.02%  "
Set-TextValue $ws.Range("D44") "0.0₆0306"
Set-TextValue $ws.Range("E44") "  -0.38%  "
Set-TextValue $ws.Range("D45") "152.43"
Set-TextValue $ws.Range("E45") "  +0.57%  "
Set-TextValue $ws.Range("E46") "  +1.73%  "
Set-TextValue $ws.Range("E47") "  -2.23%  "
Set-TextValue $ws.Range("D48") "0.602"
Set-TextValue $ws.Range("D49") "19.88"
Set-TextValue $ws.Range("E49") "  -3.76%  "
Set-TextValue $ws.Range("E50") "  -1.96%  "
Set-TextValue $ws.Range("E51") "  -2.03%  "
